$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> (new Price text, new Volume(1h) text). $null means "leave unchanged".
$updates = @{
    2 = @("29.803.49", "  +8.45%  ")
    3 = @("1.951.96", "  +6.82%  ")
    4 = @($null, "  -0.53%  ")
    5 = @("342.34", "  +2.62%  ")
    6 = @($null, "  -0.45%  ")
    7 = @($null, "  +4.16%  ")
    8 = @("0.4154", "  +8.79%  ")
    9 = @("47.82", "  +3.59%  ")
    10 = @("0.08262", "  +5.64%  ")
    11 = @("1.039", "  +8.44%  ")
    12 = @("22.76", "  +8.47%  ")
    13 = @("1.951.79", "  +6.46%  ")
    14 = @("6.177", "  +5.94%  ")
    15 = @("7.400", "  +5.03%  ")
    16 = @("92.10", "  +2.94%  ")
    17 = @($null, "  -0.47%  ")
    18 = @("0.00001059", "  +3.96%  ")
    19 = @("0.06702", "  +1.72%  ")
    20 = @("18.04", "  +5.75%  ")
    21 = @("1.001", "  -0.41%  ")
    22 = @("29.765.47", "  +8.33%  ")
    23 = @($null, "  +5.76%  ")
    24 = @($null, "  +4.78%  ")
    25 = @("2.282", "  +0.15%  ")
    26 = @("2.183.49", "  +5.36%  ")
    27 = @("161.80", "  +1.79%  ")
    28 = @("20.26", "  +4.93%  ")
    29 = @("2.180", "  +7.42%  ")
    30 = @("5.702", "  +8.22%  ")
    31 = @("122.84", "  +4.42%  ")
    32 = @($null, "  +8.94%  ")
    33 = @("0.09650", "  +2.89%  ")
    34 = @("1.480", "  +12.97%  ")
    35 = @("3.683", "  +2.93%  ")
    36 = @("5.524", "  +6.46%  ")
    37 = @("0.06289", "  +5.93%  ")
    38 = @("0.02319", "  +6.67%  ")
    39 = @("8.489", "  +4.77%  ")
    40 = @("1.189", "  +4.18%  ")
    41 = @("0.6109", "  +6.99%  ")
    42 = @("10.75", "  +8.81%  ")
    43 = @($null, "  +4.30%  ")
    44 = @("1.001", "  -0.39%  ")
    45 = @("2.407", "  +36.01%  ")
    46 = @("1.269", "  -0.08%  ")
    47 = @("0.5726", "  +6.62%  ")
    48 = @("12.49", "  +6.04%  ")
    49 = @("1.989", "  +5.01%  ")
    50 = @("0.07368", "  +7.70%  ")
    51 = @("113.42", "  +3.23%  ")
}

foreach ($row in $updates.Keys) {
    $pair = $updates[$row]
    $newPrice = $pair[0]
    $newVolume = $pair[1]
    if ($null -ne $newPrice) {
        # Force text so values like "342.34" or "1.951.96" are not reinterpreted as numbers/dates.
        $priceCell = $ws.Cells.Item($row, 4)
        $priceCell.NumberFormat = "@"
        $priceCell.Value = $newPrice
        $priceCell.Style = "Normal"
    }
    if ($null -ne $newVolume) {
        $ws.Cells.Item($row, 5).Value = $newVolume
    }
}
